# 5th test case added to BOA app
#
# - Remove the "expedia" sheet (first tab).
# - The existing "boaSignup" sheet becomes the first tab.
# - Add a new "loginNegativeTest" sheet as the last (now active) tab,
#   containing a small negative-login test data table.

$wb = $excel.ActiveWorkbook

# Drop the old "expedia" sheet entirely.
$wb.Worksheets.Item("expedia").Delete()

# Add the new sheet after the last existing sheet (so it ends up after
# "boaSignup"), then name it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "loginNegativeTest"

# Header row.
$ws.Cells.Item(1, 1).Value = "OnlineId"
$ws.Cells.Item(1, 2).Value = "Password"

# Column A (OnlineId) data.
$ws.Cells.Item(2, 1).Value = "ABCDE"
$ws.Cells.Item(3, 1).Value = "FGHIJK"
$ws.Cells.Item(4, 1).Value = "LMNOPQ"

# Column B (Password) data.
$ws.Cells.Item(2, 2).Value = "abcd123#"
$ws.Cells.Item(3, 2).Value = "FHG1234$"
$ws.Cells.Item(4, 2).Value = "pqrst123#"

# Column C (ErrContains) header + data.
$ws.Cells.Item(1, 3).Value = "ErrContains"
$ws.Cells.Item(2, 3).Value = "does not match"
$ws.Cells.Item(3, 3).Value = "does not match"
$ws.Cells.Item(4, 3).Value = "does not match"

# Match the selection left on the new active sheet.
$ws.Range("D12").Select()
